# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets carry the same underlying data, so the same row/value updates
# are applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    5  = 79
    7  = 1222
    8  = 1509
    9  = 333
    12 = 139
    18 = 317
    19 = 1706
    23 = 655
    26 = 4113
    29 = 256
    30 = 1071
    33 = 468
    35 = 216
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
